# Update vm_pu results for the 380 kV case (rows 2-25, columns B-F, I-N).
# Column order matches the original sheet layout; column H has no data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

$data = @{
    2 = @(1.02, 1.015084137785888, 1.022087458887605, 1.016789755207229, 1.031379151959334, 1.027856425568314, 1.020311677675691, 1.024922873412079, 1.019640877351437, 1.034187440470943, 1.010885623835997)
    3 = @(1.02, 1.015908660707405, 1.022701767512549, 1.017485464114561, 1.033484698000304, 1.028024443252112, 1.020771720121332, 1.025344299351542, 1.02014234503132, 1.036098064794355, 1.011039298832621)
    4 = @(1.02, 1.016442678424688, 1.023099366411004, 1.01793646191169, 1.03484167562757, 1.028131507387049, 1.021069264025487, 1.025616399319271, 1.020466982454274, 1.037328573142726, 1.011138646340427)
    5 = @(1.02, 1.016667297394194, 1.023266540252471, 1.018126258334186, 1.035410873689784, 1.028176122095616, 1.02119431877792, 1.025730648464343, 1.020603496222721, 1.037844519566617, 1.011180390138239)
    6 = @(1.02, 1.016705018796649, 1.023294610819279, 1.01815813750316, 1.035506370619909, 1.028183589965306, 1.021215314087121, 1.025749823082825, 1.020626419585812, 1.037931070328595, 1.01118739781488)
    7 = @(1.02, 1.01644567933213, 1.02310160010558, 1.017938997208419, 1.034849286254934, 1.028132105082004, 1.021070935142433, 1.025617926478998, 1.020468806416665, 1.037335472551563, 1.011139204209425)
    8 = @(1.02, 1.015362685887752, 1.022295046208884, 1.017024701219309, 1.032091877889973, 1.027913551152555, 1.020467178660187, 1.025065418468549, 1.019810318308593, 1.034834361212831, 1.010937577619989)
    9 = @(1.02, 1.0134581519896, 1.020874592173676, 1.015419984945692, 1.027189872362439, 1.02751571120355, 1.019402277278631, 1.024087313516823, 1.018651192055721, 1.030381460319251, 1.010581602672693)
    10 = @(1.02, 1.012191093028334, 1.019928196252319, 1.014354538316065, 1.023890989568029, 1.02724186259601, 1.018691696727212, 1.023432217713777, 1.017879299088782, 1.027380417523079, 1.010343839476615)
    11 = @(1.02, 1.011643076034334, 1.019518539400238, 1.01389423712305, 1.022454784957048, 1.027121221850237, 1.018383859840849, 1.023147837486483, 1.017545272656863, 1.026072840084491, 1.01024078195399)
    12 = @(1.02, 1.011439613162708, 1.019366395943841, 1.013723418725557, 1.021920110616982, 1.02707609916055, 1.018269493247857, 1.023042097922, 1.017421232282377, 1.025585895174092, 1.010202486283023)
    13 = @(1.02, 1.011483252327393, 1.019399030258666, 1.013760052686529, 1.022034855164771, 1.027085792243821, 1.018294026271795, 1.023064784308926, 1.017447837909807, 1.025690403746883, 1.010210701533859)
    14 = @(1.02, 1.011626255787072, 1.019505962736695, 1.013880113992518, 1.022410613387823, 1.027117498354668, 1.01837440671463, 1.023139099226118, 1.017535018786247, 1.026032614817688, 1.010237616737433)
    15 = @(1.02, 1.011714377552914, 1.019571850207071, 1.013954108688857, 1.022641969766212, 1.027136992234852, 1.018423928819185, 1.023184872776381, 1.017588738048775, 1.026243295348886, 1.010254198031607)
    16 = @(1.02, 1.012227476401061, 1.019955386795081, 1.014385109083321, 1.023986138947763, 1.027249825545245, 1.018712123703787, 1.02345107593598, 1.017901471765061, 1.027467023407028, 1.010350676879001)
    17 = @(1.02, 1.012549498361295, 1.020196006680199, 1.014655744579163, 1.024827195393797, 1.027320049679481, 1.018892860644194, 1.023617865542336, 1.018097697576043, 1.02823244122872, 1.010411167665015)
    18 = @(1.02, 1.012737388860058, 1.020336369553224, 1.014813702367506, 1.025317022567164, 1.027360811330537, 1.018998266867247, 1.023715081686386, 1.018212172831337, 1.028678116738073, 1.010446440823816)
    19 = @(1.02, 1.012801464910723, 1.020384231917139, 1.014867578935498, 1.025483915416099, 1.027374676307626, 1.019034205180507, 1.02374821810516, 1.018251209318264, 1.028829949346904, 1.010458466341781)
    20 = @(1.02, 1.01251494218592, 1.020170189057915, 1.014626697529996, 1.024737035511244, 1.027312535878115, 1.018873470785916, 1.023599977775662, 1.018076642324339, 1.028150400042714, 1.010404678618871)
    21 = @(1.02, 1.011584142203292, 1.01947447319437, 1.013844754561313, 1.022299995494544, 1.02710817030199, 1.018350737289262, 1.023117218319094, 1.017509345305176, 1.025931877040575, 1.010229691307794)
    22 = @(1.02, 1.010999461523146, 1.01903717269252, 1.013354030411953, 1.020760746861388, 1.026977875446112, 1.018021945344096, 1.022813062916163, 1.017152848096445, 1.02452973983805, 1.010119580139859)
    23 = @(1.02, 1.011309359428334, 1.019268982086171, 1.013614085535773, 1.021577406341652, 1.027047118536738, 1.018196256242976, 1.02297436073484, 1.017341816279218, 1.025273739518401, 1.010177960597297)
    24 = @(1.02, 1.012530556447568, 1.020181854889207, 1.014639822330542, 1.024777777187538, 1.027315931656165, 1.018882232273585, 1.023608060702198, 1.018086156224617, 1.028187473331367, 1.0104076107708)
    25 = @(1.02, 1.01395005959941, 1.021241715175568, 1.015834078254979, 1.028462459674692, 1.027620076819937, 1.019677695838351, 1.02434071106201, 1.018950705770512, 1.031538228731713, 1.010673710352054)
}

foreach ($row in $data.Keys) {
    $rowValues = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $row).Value = $rowValues[$i]
    }
}
